# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The workbook stores one row per game starting at row 2 (row 1 is the
# header: date, TB, PC, dS0, dSF, K, IP, I0, IF). Column G ("K") previously
# held a different ("Strike#") statistic; this recomputes/overwrites it with
# the correct K values for every game row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values for rows 2..70 (row index => value), keyed by worksheet row.
$kValues = @{
    2  = 1
    3  = 0
    4  = 4
    5  = 1
    6  = 1
    7  = 0
    8  = 1
    9  = 2
    10 = 1
    11 = 1
    12 = 3
    13 = 0
    14 = 0
    15 = 1
    16 = 2
    17 = 1
    18 = 0
    19 = 0
    20 = 1
    21 = 0
    22 = 1
    23 = 0
    24 = 1
    25 = 2
    26 = 0
    27 = 0
    28 = 0
    29 = 1
    30 = 0
    31 = 4
    32 = 0
    33 = 3
    34 = 3
    35 = 1
    36 = 1
    37 = 1
    38 = 1
    39 = 1
    40 = 2
    41 = 2
    42 = 1
    43 = 1
    44 = 2
    45 = 1
    46 = 1
    47 = 1
    48 = 2
    49 = 1
    50 = 1
    51 = 1
    52 = 2
    53 = 1
    54 = 0
    55 = 1
    56 = 0
    57 = 1
    58 = 1
    59 = 1
    60 = 2
    61 = 1
    62 = 1
    63 = 2
    64 = 2
    65 = 4
    66 = 1
    67 = 2
    68 = 0
    69 = 3
    70 = 3
}

foreach ($r in ($kValues.Keys | Sort-Object)) {
    $ws.Cells.Item($r, 7).Value = $kValues[$r]
}
